$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 7 de Abril de 2020 a las 14:52"

# Row 4
$ws.Range("B4").Value = 367719
$ws.Range("C4").Value = 715
$ws.Range("E4").Value = 336962

# Row 7
$ws.Range("E7").Value = 65814
$ws.Range("G7").Value = 12
$ws.Range("H7").Value = 1822

# Row 27
$ws.Range("F27").Value = 127

# Row 39
$ws.Range("B39").Value = 2795
$ws.Range("C39").Value = 190
$ws.Range("D39").Value = 615
$ws.Range("E39").Value = 2139
$ws.Range("G39").Value = 3
$ws.Range("H39").Value = 41

# Row 52
$ws.Range("D52").Value = 338
$ws.Range("E52").Value = 1237

# Row 63
$ws.Range("A63").Value = "Eslovenia"
$ws.Range("B63").Value = 1059
$ws.Range("C63").Value = 38
$ws.Range("D63").Value = 102
$ws.Range("E63").Value = 927
$ws.Range("F63").Value = 114
$ws.Range("H63").Value = 30

# Row 64
$ws.Range("A64").Value = "Irak"
$ws.Range("B64").Value = 1031
$ws.Range("D64").Value = 344
$ws.Range("E64").Value = 623
$ws.Range("F64").Value = 0
$ws.Range("H64").Value = 64

# Row 78
$ws.Range("A78").Value = "Republica de Macedonia"
$ws.Range("B78").Value = 599
$ws.Range("C78").Value = 29
$ws.Range("D78").Value = 30
$ws.Range("E78").Value = 543
$ws.Range("F78").Value = 15
$ws.Range("G78").Value = 3
$ws.Range("H78").Value = 26

# Row 79
$ws.Range("A79").Value = "Tunez"
$ws.Range("B79").Value = 596
$ws.Range("C79").Value = 0
$ws.Range("D79").Value = 5
$ws.Range("E79").Value = 569
$ws.Range("F79").Value = 39
$ws.Range("H79").Value = 22

# Row 80
$ws.Range("A80").Value = "Eslovaquia"
$ws.Range("B80").Value = 581
$ws.Range("C80").Value = 47
$ws.Range("D80").Value = 9
$ws.Range("E80").Value = 570
$ws.Range("F80").Value = 3
$ws.Range("H80").Value = 2

# Row 114
$ws.Range("A114").Value = "Kenia"
$ws.Range("B114").Value = 172
$ws.Range("C114").Value = 14
$ws.Range("D114").Value = 4
$ws.Range("E114").Value = 162
$ws.Range("F114").Value = 2
$ws.Range("H114").Value = 6

# Row 115
$ws.Range("A115").Value = "Venezuela"
$ws.Range("B115").Value = 165
$ws.Range("D115").Value = 65
$ws.Range("E115").Value = 93
$ws.Range("F115").Value = 6
$ws.Range("H115").Value = 7

# Row 116
$ws.Range("A116").Value = "Mayotte"
$ws.Range("C116").Value = 0
$ws.Range("D116").Value = 15
$ws.Range("E116").Value = 147
$ws.Range("F116").Value = 3
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 2

# Row 117
$ws.Range("A117").Value = "Banglades"
$ws.Range("B117").Value = 164
$ws.Range("C117").Value = 41
$ws.Range("D117").Value = 33
$ws.Range("E117").Value = 114
$ws.Range("F117").Value = 1
$ws.Range("G117").Value = 5
$ws.Range("H117").Value = 17

# Row 118
$ws.Range("A118").Value = "Consejo Danes para los Refugiados"
$ws.Range("B118").Value = 161
$ws.Range("D118").Value = 5
$ws.Range("E118").Value = 138
$ws.Range("F118").Value = 0
$ws.Range("H118").Value = 18
